# Applies the "Added functions getRoot() and findNext() ..." journal entry
# to the Jesse worksheet (row 10), matching the commit's xlsx diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jesse")

$noteText = "Added functions getRoot() and findNext() to Building.h to premit traversal. `nAdded data members, inheretance, and Deploy() function to States.h to handle traversal/generation. `nAdded comments for coding suggestions in States.h"

# Write the new row's values FIRST (while the cells are still blank/default
# formatted), so the SUM(B4:B200)/(.../60) dependency chain picks the row up
# immediately. Formatting the blank cells before giving them a value can
# leave the range aggregation stale even after a recalculation.
$ws.Cells.Item(10, 1).Value = 43074
$ws.Cells.Item(10, 2).Value = 180
$ws.Cells.Item(10, 3).Value = $noteText

# Now copy the per-cell formatting from row 9 into row 10.
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(10).RowHeight = 99.75

# Update the view state to match the committed selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("F10").Select()

$excel.CalculateFullRebuild()
